# Auto-generated Word COM-interop script implementing the diff.
$d = $word.ActiveDocument

# Colección ESTADÍSTICAS final paragraph - split runs with gramStart/End
$xml_32 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">No definido </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>todavía :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>’c</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(32).Range.InsertXML($xml_32)

# fecha_creacion - add proofErr
$xml_27 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="312" w:lineRule="atLeast"/><w:rPr><w:rStyle w:val="red-ui-debug-msg-element"/><w:color w:val="333333"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-object-key"/><w:color w:val="792E90"/></w:rPr><w:t>fecha_creacion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-tools"/><w:color w:val="333333"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rStyle w:val="apple-converted-space"/><w:color w:val="333333"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-type-string"/><w:color w:val="B72828"/></w:rPr><w:t>"2021-12-28T11:45:22.103Z"</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(27).Range.InsertXML($xml_27)

# id_reserva - add proofErr
$xml_26 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="312" w:lineRule="atLeast"/><w:rPr><w:rStyle w:val="red-ui-debug-msg-element"/><w:color w:val="333333"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-object-key"/><w:color w:val="792E90"/></w:rPr><w:t>id_reserva</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-tools"/><w:color w:val="333333"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rStyle w:val="apple-converted-space"/><w:color w:val="333333"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-type-number"/><w:color w:val="2033D6"/></w:rPr><w:t>1</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(26).Range.InsertXML($xml_26)

# hora_fin - add proofErr
$xml_25 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="312" w:lineRule="atLeast"/><w:rPr><w:rStyle w:val="red-ui-debug-msg-element"/><w:color w:val="333333"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-object-key"/><w:color w:val="792E90"/></w:rPr><w:t>hora_fin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-tools"/><w:color w:val="333333"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rStyle w:val="apple-converted-space"/><w:color w:val="333333"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-type-number"/><w:color w:val="2033D6"/></w:rPr><w:t>18</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(25).Range.InsertXML($xml_25)

# hora_inicio - add proofErr
$xml_24 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="312" w:lineRule="atLeast"/><w:rPr><w:rStyle w:val="red-ui-debug-msg-element"/><w:color w:val="333333"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-object-key"/><w:color w:val="792E90"/></w:rPr><w:t>hora_inicio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-tools"/><w:color w:val="333333"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rStyle w:val="apple-converted-space"/><w:color w:val="333333"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-type-number"/><w:color w:val="2033D6"/></w:rPr><w:t>15</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(24).Range.InsertXML($xml_24)

# id_telegram - add proofErr
$xml_21 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="312" w:lineRule="atLeast"/><w:rPr><w:rStyle w:val="red-ui-debug-msg-element"/><w:color w:val="333333"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-object-key"/><w:color w:val="792E90"/></w:rPr><w:t>id_telegram</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-tools"/><w:color w:val="333333"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rStyle w:val="apple-converted-space"/><w:color w:val="333333"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rStyle w:val="red-ui-debug-msg-type-number"/><w:color w:val="2033D6"/></w:rPr><w:t>123</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(21).Range.InsertXML($xml_21)

# date -> fecha_inclusion
$xml_14 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="312" w:lineRule="atLeast"/><w:rPr><w:color w:val="333333"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="792E90"/></w:rPr><w:t>fecha</w:t></w:r><w:r><w:rPr><w:color w:val="792E90"/></w:rPr><w:t>_inclusion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="333333"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:color w:val="B72828"/></w:rPr><w:t>"2021-12-30T13:33:36.524Z"</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(14).Range.InsertXML($xml_14)

# Suscipcion -> suscripcion, Arranque -> Maker
$xml_13 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="312" w:lineRule="atLeast"/><w:rPr><w:color w:val="333333"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="792E90"/></w:rPr><w:t>suscripcion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="333333"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:color w:val="B72828"/></w:rPr><w:t>"</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="B72828"/></w:rPr><w:t>Maker</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="B72828"/></w:rPr><w:t>"</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(13).Range.InsertXML($xml_13)

# Acceso CNC..ID NFC block restructure (paragraphs 9-12 -> 4 new paragraphs)
$rng9 = $d.Range($d.Paragraphs(9).Range.Start, $d.Paragraphs(12).Range.End)
$xml_9 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="312" w:lineRule="atLeast"/><w:rPr><w:color w:val="333333"/><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="792E90"/><w:lang w:val="it-IT"/></w:rPr><w:t>acceso_CNC</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="333333"/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:color w:val="2033D6"/><w:lang w:val="it-IT"/></w:rPr><w:t>false</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="312" w:lineRule="atLeast"/><w:rPr><w:color w:val="333333"/><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="792E90"/><w:lang w:val="it-IT"/></w:rPr><w:t>acceso_3D</w:t></w:r><w:r><w:rPr><w:color w:val="333333"/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="2033D6"/><w:lang w:val="it-IT"/></w:rPr><w:t>true</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:spacing w:line="312" w:lineRule="atLeast"/><w:rPr><w:color w:val="333333"/><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="792E90"/><w:lang w:val="it-IT"/></w:rPr><w:t>password</w:t></w:r><w:r><w:rPr><w:color w:val="333333"/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:color w:val="B72828"/><w:lang w:val="it-IT"/></w:rPr><w:t>"</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="B72828"/><w:lang w:val="it-IT"/></w:rPr><w:t>pato</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="B72828"/><w:lang w:val="it-IT"/></w:rPr><w:t>"</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="312" w:lineRule="atLeast"/><w:rPr><w:color w:val="333333"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="792E90"/></w:rPr><w:t>id_NFC</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="333333"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="666666"/></w:rPr><w:t>null</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng9.InsertXML($xml_9)

# Fecha Vencimiento -> fecha_vencimiento, date value change
$xml_8 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="312" w:lineRule="atLeast"/><w:rPr><w:color w:val="333333"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="792E90"/></w:rPr><w:t>fecha_vencimiento</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="333333"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:color w:val="B72828"/></w:rPr><w:t>"2022-01-30T23:00:00.000Z"</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(8).Range.InsertXML($xml_8)

# Nombre -> nombre, maria -> aguacate
$xml_7 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:color w:val="792E90"/></w:rPr><w:t>nombre</w:t></w:r><w:r><w:rPr><w:color w:val="333333"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:color w:val="B72828"/></w:rPr><w:t>"aguacate"</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(7).Range.InsertXML($xml_7)
